$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (rows 2-17, columns B:G) after filtering save games.
$data = @(
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455),
    @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455),
    @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.964545797025059),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027),
    @(1.445647641019636, 1.626987699542094, 3993.344853322108, 13.86384647080068, 1, 4010.28133513347),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248),
    @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 17.65757632934944),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 19.48425592650926),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027),
    @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248),
    @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 21.98653043760045),
    @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
